# Weekly update: insert a new price observation row for
# "Vega Monumental Concepción - Mango" as row 95, pushing the
# previously existing rows 95-98 down to 96-99.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 95 (shifts rows 95:98 -> 96:99,
# and extends the used range/dimension to A1:T99).
$ws.Rows.Item(95).Insert()

# Populate the new row 95 with this week's record.
$ws.Range("A95").Value = 11
$ws.Range("B95").Value = "Vega Monumental Concepción"
$ws.Range("C95").Value = "Bíobío"
$ws.Range("D95").Value = 44610
$ws.Range("E95").Value = 8
$ws.Range("F95").Value = "Fruta"
$ws.Range("G95").Value = 100108
$ws.Range("H95").Value = "Tropicales y subtropicales"
$ws.Range("I95").Value = 100108002
$ws.Range("J95").Value = "Mango"
$ws.Range("K95").Value = "Sin especificar"
$ws.Range("L95").Value = "Primera"
$ws.Range("M95").Value = 180
$ws.Range("N95").Value = 7500
$ws.Range("O95").Value = 8000
$ws.Range("P95").Value = 7722
$ws.Range("Q95").Value = "$/bandeja 4 kilos"
$ws.Range("R95").Value = "Perú"
$ws.Range("S95").Value = 1930
$ws.Range("T95").Value = 4
